$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet: swap the SFORZA / TORENBEEK_1982 rows (labels + values) ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")

$wsFuselage.Range("A23").Value = "TORENBEEK_1982"
$wsFuselage.Range("A24").Value = "SFORZA"

$wsFuselage.Range("C23").Value = 10.594739999999998
$wsFuselage.Range("C24").Value = 11.500334910927485

# --- WING sheet: swap the SFORZA / TORENBEEK_1982 rows (labels + values) ---
$wsWing = $wb.Worksheets.Item("WING")

$wsWing.Range("A23").Value = "TORENBEEK_1982"
$wsWing.Range("A24").Value = "SFORZA"

$wsWing.Range("C23").Value = 1.133712717373045
$wsWing.Range("C24").Value = 0.8092048979331106

$wsWing.Range("A27").Value = "TORENBEEK_1982"
$wsWing.Range("A28").Value = "SFORZA"

$wsWing.Range("C27").Value = 5.087959999999999
$wsWing.Range("C28").Value = 10.850510037388545
